$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column headers in row 1 so that the "_old"/"_new" suffixes used to
# distinguish the two compared format versions are replaced by the concrete
# format-version identifiers of the input files ("FV2404" / "FV2410").
$lastCol = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2404")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2410")
        }
    }
}

# Turn the used range into a real Excel Table ("ListObject") with headers,
# matching the freshly renamed column captions.
$lastRow = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
